$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.823.30'
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").Value = '3.502.58'
$ws.Range("E3").Value = '  -1.61%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '''608.01'
$ws.Range("E5").Value = '  +3.80%  '
$ws.Range("D6").Value = '''191.59'
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("D7").Value = '''0.626'
$ws.Range("E7").Value = '  +0.69%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = '''0.213'
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("D10").Value = '''0.663'
$ws.Range("E10").Value = '  +3.02%  '
$ws.Range("D11").Value = '''53.47'
$ws.Range("E11").Value = '  -0.99%  '
$ws.Range("D12").Value = '''0.0000307'
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("D13").Value = '''9.60'
$ws.Range("E13").Value = '  +2.22%  '
$ws.Range("D14").Value = '4.065.61'
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").Value = '''622.87'
$ws.Range("E15").Value = '  +9.94%  '
$ws.Range("D16").Value = '69.901.85'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '''12.70'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").Value = '''18.88'
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("D19").Value = '3.511.08'
$ws.Range("E19").Value = '  -1.34%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '''0.992'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = '''17.76'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = '''105.94'
$ws.Range("E23").Value = '  +12.81%  '
$ws.Range("D24").Value = '''4.65'
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("E25").Value = '  +2.98%  '
$ws.Range("E26").Value = '  +5.00%  '
$ws.Range("D27").Value = '''10.99'
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").Value = '''9.80'
$ws.Range("E28").Value = '  +5.45%  '
$ws.Range("D29").Value = '''34.26'
$ws.Range("E29").Value = '  +5.74%  '
$ws.Range("D30").Value = '''7.07'
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("D31").Value = '''12.54'
$ws.Range("E31").Value = '  +2.99%  '
$ws.Range("D32").Value = '''4.12'
$ws.Range("E32").Value = '  +5.15%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").Value = '''64.03'
$ws.Range("E34").Value = '  +1.60%  '
$ws.Range("D35").Value = '3.718.27'
$ws.Range("E35").Value = '  +2.21%  '
$ws.Range("D36").Value = '''3.09'
$ws.Range("E36").Value = '  -4.57%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").Value = '''517.53'
$ws.Range("E38").Value = '  -1.89%  '
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '0.0₃0793'
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = '''0.390'
$ws.Range("E40").Value = '  -3.86%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''36.73'
$ws.Range("E41").Value = '  -3.71%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''3.58'
$ws.Range("E42").Value = '  +1.26%  '
$ws.Range("E43").Value = '  -0.70%  '
$ws.Range("D44").Value = '''0.0463'
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("D45").Value = '''2.86'
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("E46").Value = '  +2.53%  '
$ws.Range("D47").Value = '''3.32'
$ws.Range("E47").Value = '  -3.94%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").Value = '''1.00'
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '''8.75'
$ws.Range("E49").Value = '  -4.88%  '
$ws.Range("D50").Value = '''132.23'
$ws.Range("E50").Value = '  -1.49%  '
$ws.Range("D51").Value = '''1.35'
$ws.Range("E51").Value = '  -6.61%  '
